$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.260.29"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").Value = "3.415.01"
$ws.Range("E3").Value = "  +0.86%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.48"
$ws.Range("E5").Value = "  -1.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.08"
$ws.Range("E6").Value = "  -2.91%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "3.408.67"
$ws.Range("E8").Value = "  +0.75%  "

$ws.Range("E9").Value = "  -0.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.196"
$ws.Range("E10").Value = "  +0.27%  "

$ws.Range("E11").Value = "  -1.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.73"
$ws.Range("E12").Value = "  -0.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000278"
$ws.Range("E13").Value = "  -1.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "693.07"
$ws.Range("E14").Value = "  +0.94%  "

$ws.Range("D15").Value = "3.963.38"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.61"
$ws.Range("E16").Value = "  -0.26%  "

$ws.Range("D17").Value = "69.322.43"
$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").Value = "3.422.37"
$ws.Range("E18").Value = "  +1.56%  "

$ws.Range("E19").Value = "  +0.71%  "

$ws.Range("E20").Value = "  -0.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.34"
$ws.Range("E21").Value = "  -0.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.894"
$ws.Range("E22").Value = "  -0.61%  "

$ws.Range("E23").Value = "  +0.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.87"
$ws.Range("E24").Value = "  -1.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "100.61"
$ws.Range("E25").Value = "  -4.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.87"

$ws.Range("E27").Value = "  -2.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.54"
$ws.Range("E28").Value = "  -0.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.30"
$ws.Range("E29").Value = "  -3.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.71"
$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.94"
$ws.Range("E31").Value = "  -1.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "567.85"
$ws.Range("E32").Value = "  +1.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.66"
$ws.Range("E33").Value = "  -0.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.96"
$ws.Range("E34").Value = "  -2.17%  "

$ws.Range("E35").Value = "  -2.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.10"
$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").Value = "3.577.29"
$ws.Range("E38").Value = "  -4.00%  "

$ws.Range("E39").Value = "  -2.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.74"
$ws.Range("E40").Value = "  -1.07%  "

$ws.Range("D41").Value = "0.0₃0722"
$ws.Range("E41").Value = "  +2.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.25"
$ws.Range("E42").Value = "  -0.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.64"
$ws.Range("E43").Value = "  -1.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.330"
$ws.Range("E44").Value = "  -2.99%  "

$ws.Range("E45").Value = "  -0.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.44"
$ws.Range("E46").Value = "  +2.55%  "

$ws.Range("E47").Value = "  -1.58%  "

$ws.Range("E48").Value = "  -1.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.58"
$ws.Range("E50").Value = "  -0.69%  "

$ws.Range("E51").Value = "  +1.24%  "
